$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Type" column data/header (column G: "Type" header, "Pool"
# values) without shifting later columns.
$ws.Range("G1:G3").ClearContents()

# The now-unused "Folio No" column (column H) is removed entirely,
# shifting Instrument/Currency/Sector/Investment Domicile/Custom Field
# columns one position to the left.
$ws.Range("H1").EntireColumn.Delete()

$ws.Range("G:G").Select()
